# Fix the "Department" data-validation column:
#  - it used to be a hard-coded inline list (duplicated many times because of a
#    bug that kept re-inserting the same validations)
#  - now it should read its list of choices off a new hidden "departments"
#    sheet, which also carries two extra departments that weren't in the old
#    inline list.
# The old "Confirm Passsword" column (J) is also removed - the "Department"
# column shifts from K into J as a result.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the new "departments" lookup sheet right after user_department, and
#    populate it with the master list of departments.
# ---------------------------------------------------------------------------
$deptSheet = $wb.Worksheets.Add($null, $ws1)
$deptSheet.Name = "departments"

$deptValues = @(
    "Department",
    "Department of Commerce",
    "Department of Business Studies",
    "Department of Cultural Studies",
    "Department of Media Communications",
    "Department of Arts",
    "Department of Mathematics",
    "Department of Physical Education",
    "Department of Computer ScienceDepartment of Computer ScienceDepartment of Computer ScienceDepartment of Computer ScienceDepartment of Computer ScienceDepartment of Computer Science"
)
for ($i = 0; $i -lt $deptValues.Length; $i++) {
    $deptSheet.Cells.Item($i + 1, 1).Value = $deptValues[$i]
}
$deptSheet.Columns.Item(1).ColumnWidth = 38.6
$deptSheet.Range("A2").Select() | Out-Null

# The lookup sheet is an implementation detail - keep it hidden from users.
$deptSheet.Visible = $false

# ---------------------------------------------------------------------------
# 2. Remove the old "Confirm Passsword" column (J). This shifts the old
#    "Department" column (K) left into J, carrying its header/width with it.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(10).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. Rebuild the data validations cleanly: Sex (F), Salutation (G) and Health
#    (H) rules stay the same; the Department rule (now on J) switches from a
#    hard-coded inline list to a reference into the 'departments' sheet.
# ---------------------------------------------------------------------------
$ws1.Cells.Validation.Delete() | Out-Null

$ws1.Range("F2:F100").Validation.Add(3, 1, 1, '"Male,Female,Unspecified"')
$ws1.Range("F2:F100").Validation.IgnoreBlank = $false

$ws1.Range("G2:G100").Validation.Add(3, 1, 1, '"Dr,Mr,Ms,Prof,Rev"')
$ws1.Range("G2:G100").Validation.IgnoreBlank = $false

$ws1.Range("H2:H100").Validation.Add(3, 1, 1, '"Fine,Not Fine"')
$ws1.Range("H2:H100").Validation.IgnoreBlank = $false

$ws1.Range("J2:J100").Validation.Add(3, 1, 1, "='departments'!`$A`$2:`$A`$9")
$ws1.Range("J2:J100").Validation.IgnoreBlank = $true

# ---------------------------------------------------------------------------
# 4. Restore the expected selection state on the main sheet (whole column J).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(10).Select() | Out-Null
$ws1.Activate() | Out-Null
